$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on the numeric-looking cells before writing their
# new values, so Excel stores them as text (matching the workbooks original
# inline-string representation) instead of converting them to numbers.
$numericCells = @(
    'D2',
    'G2',
    'D3',
    'G3',
    'D4',
    'G4',
    'D5',
    'G5',
    'D6',
    'G6',
    'G7',
    'D8',
    'G8',
    'D9',
    'G9',
    'D10',
    'G10',
    'D11',
    'G11',
    'D12',
    'G12',
    'D13',
    'G13',
    'D14',
    'G14',
    'D15',
    'G15',
    'D16',
    'G16',
    'D17',
    'G17',
    'D18',
    'G18',
    'D19',
    'G19',
    'D20',
    'G20',
    'D21',
    'G21',
    'D22',
    'G22',
    'D23',
    'G23',
    'D24',
    'G24',
    'G25',
    'G26',
    'G27',
    'G28',
    'G29',
    'G30',
    'G31',
    'G32',
    'G33',
    'G34',
    'G35',
    'G36',
    'G37',
    'G38',
    'G39',
    'D40',
    'G40',
    'D41',
    'G41',
    'D42',
    'G42',
    'G43',
    'D44',
    'G44',
    'D45',
    'G45',
    'G46',
    'G47',
    'G48',
    'G49',
    'G50',
    'G51'
)
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '244.79'
$ws.Range('G2').Value = '9'
$ws.Range('D3').Value = '23.08'
$ws.Range('G3').Value = '9'
$ws.Range('D4').Value = '5.415'
$ws.Range('G4').Value = '9'
$ws.Range('D5').Value = '0.05993'
$ws.Range('G5').Value = '9'
$ws.Range('D6').Value = '3.395'
$ws.Range('G6').Value = '9'
$ws.Range('G7').Value = '9'
$ws.Range('D8').Value = '0.9271'
$ws.Range('G8').Value = '9'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').Value = '0.1435'
$ws.Range('E9').Value = '8WazirXWRX'
$ws.Range('G9').Value = '9'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').Value = '0.07441'
$ws.Range('E10').Value = '9MandalaExchangeTokenMDX'
$ws.Range('G10').Value = '9'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '0.03377'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('G11').Value = '9'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.03037'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('G12').Value = '9'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.09354'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('G13').Value = '9'
$ws.Range('B14').Value = 'MCDex'
$ws.Range('C14').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D14').Value = '3.939'
$ws.Range('E14').Value = '13MCDexMCB'
$ws.Range('G14').Value = '9'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001593'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('G15').Value = '9'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').Value = '0.04812'
$ws.Range('E16').Value = '15CoinExTokenCET'
$ws.Range('G16').Value = '9'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').Value = '0.0005942'
$ws.Range('E17').Value = '16OneONEWorstin24h'
$ws.Range('G17').Value = '9'
$ws.Range('D18').Value = '0.005628'
$ws.Range('G18').Value = '9'
$ws.Range('D19').Value = '0.004159'
$ws.Range('G19').Value = '9'
$ws.Range('D20').Value = '0.0009833'
$ws.Range('G20').Value = '9'
$ws.Range('D21').Value = '0.00007702'
$ws.Range('G21').Value = '9'
$ws.Range('D22').Value = '3.659'
$ws.Range('G22').Value = '9'
$ws.Range('D23').Value = '6.450'
$ws.Range('G23').Value = '9'
$ws.Range('D24').Value = '2.186'
$ws.Range('G24').Value = '9'
$ws.Range('G25').Value = '9'
$ws.Range('G26').Value = '9'
$ws.Range('G27').Value = '9'
$ws.Range('G28').Value = '9'
$ws.Range('G29').Value = '9'
$ws.Range('G30').Value = '9'
$ws.Range('G31').Value = '9'
$ws.Range('G32').Value = '9'
$ws.Range('G33').Value = '9'
$ws.Range('G34').Value = '9'
$ws.Range('G35').Value = '9'
$ws.Range('G36').Value = '9'
$ws.Range('G37').Value = '9'
$ws.Range('G38').Value = '9'
$ws.Range('G39').Value = '9'
$ws.Range('D40').Value = '0.03955'
$ws.Range('G40').Value = '9'
$ws.Range('D41').Value = '0.006222'
$ws.Range('G41').Value = '9'
$ws.Range('D42').Value = '0.1075'
$ws.Range('G42').Value = '9'
$ws.Range('G43').Value = '9'
$ws.Range('D44').Value = '0.007486'
$ws.Range('G44').Value = '9'
$ws.Range('D45').Value = '0.00005209'
$ws.Range('G45').Value = '9'
$ws.Range('G46').Value = '9'
$ws.Range('G47').Value = '9'
$ws.Range('E48').Value = '47CoinbaseStockTokenCOINBestin24h'
$ws.Range('G48').Value = '9'
$ws.Range('E49').Value = '48BOLOBOLO'
$ws.Range('G49').Value = '9'
$ws.Range('G50').Value = '9'
$ws.Range('G51').Value = '9'
